$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

# Remove IRA-era "ban new power plants" flags (set back to 0) for
# hard coal (row 2) and lignite (row 14), columns I (2028) through AE (2050).
$ws.Range("I2:AE2").Value = 0
$ws.Range("I14:AE14").Value = 0

# Update the selection to match the final saved state (H14:AE14).
$ws.Range("H14:AE14").Select()
